$d = $word.ActiveDocument

# The job-title cell in the "WORKEXPERIENCE" table currently reads
# "Jr.CreditManager". Replace it with " CreditAssociate" (i.e. the role
# title changes from "Jr. Credit Manager" to "Credit Associate").
$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute("Jr.CreditManager", $true, $false, $false, $false, $false, $true, 1, $false, " CreditAssociate", 2) | Out-Null
